$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style index 9 look) from A24 onto the new data rows A25:G29
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A25:G29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 25
$ws.Cells.Item(25, 1).Value = 43902
$ws.Cells.Item(25, 2).Value = '17:00-19:00'
$ws.Cells.Item(25, 4).Value = 'Another three Key expert practices, three advanced topics, fruitful attitudes,  enduring principles'
$ws.Cells.Item(25, 5).Value = 'Get familiar with the topics whch will be useful in the future'
$ws.Cells.Item(25, 6).Value = 'Attitude is important according to my experience. Healthy attitudes can affect your behavior and learn more from others.'
$ws.Cells.Item(25, 7).Value = 'The last course is shorter but still inspiring. Thanks for the effort of Andre and Kaj. Hope I can make full use of all the things I learned from the course to make more achievements.'
$ws.Rows.Item(25).RowHeight = 77.5

# Row 26
$ws.Cells.Item(26, 1).Value = 43903
$ws.Cells.Item(26, 2).Value = '11:00-12:00'
$ws.Cells.Item(26, 3).Value = 'Guowei Li, Dongxin Xiang'
$ws.Cells.Item(26, 4).Value = 'Figure out how to do the assignment and which part of the assignment each of us should be mainly responsible for.'
$ws.Cells.Item(26, 5).Value = 'Discuss the assignment with my partners. We discussed what we should do with the assignment and how to do it.'
$ws.Cells.Item(26, 6).Value = 'Each of us has different advantages  which can be used to do the assignment. And our discussion can inpire each other.'
$ws.Cells.Item(26, 7).Value = 'Our team started to get things done more efficiently.'
$ws.Rows.Item(26).RowHeight = 78

# Row 27
$ws.Cells.Item(27, 1).Value = 43904
$ws.Cells.Item(27, 2).Value = '10:51-12:00, 16:30-18:20'
$ws.Cells.Item(27, 4).Value = 'Find interesting test cases. Decide test cases and new test cases to use in our reports together after we found them seperately.'
$ws.Cells.Item(27, 5).Value = 'Read all the test cases in the project and find the ones I am interested in.  And decide what to use in our assignment with my partners through Wechat'
$ws.Cells.Item(27, 6).Value = 'There are not many test cases in the project, but some of the most important features were tested. I am wondering if most Android apps have not many test cases.'
$ws.Cells.Item(27, 7).Value = 'I am interested in the UI test cases. They are useful and fast. And I learned how to write them by myself.'
$ws.Rows.Item(27).RowHeight = 91

# Row 28
$ws.Cells.Item(28, 1).Value = 43905
$ws.Cells.Item(28, 2).Value = '10:00-11:00, 21:00-24:00'
$ws.Cells.Item(28, 3).Value = 'Guowei Li, Dongxin Xiang'
$ws.Cells.Item(28, 4).Value = 'Modify the reports'
$ws.Cells.Item(28, 5).Value = 'We continued finishing the reports and we shared ideas for all the reports and give suggestions to modify them.'
$ws.Cells.Item(28, 6).Value = 'We used some skills learned at Prof. Jones''s class to do the assignment. So  this assignment was little more manageable than those before. And we tried to combine knowledge in two courses to do a good job.'
$ws.Cells.Item(28, 7).Value = 'I am glad we can apply some skills in the testing course to this assignment. It gave me a chance to solve a problem by using skills in diffetent fields amd get more practice to reinforce my learning.'
$ws.Rows.Item(28).RowHeight = 150

# Row 29
$ws.Cells.Item(29, 1).Value = 43906
$ws.Cells.Item(29, 2).Value = '8:30-9:30'
$ws.Cells.Item(29, 4).Value = 'Finish the reports'
$ws.Cells.Item(29, 5).Value = 'Reviewed all the reports, check if all the suggesions we made have been written down and confirmed we all agree with all the contents.'
$ws.Cells.Item(29, 6).Value = 'One of the test case was especially interesting. I proposed an open issue last week, and I already found some clues. Then my partner found it almost impossible to be fixed by us because the developer wrote the rule in the jar file. I was glad my partner suggest we write a test case for this issue. In my opinion, It was like test-driven development. We made a pull request for the test case, so in the future when the developer of Omni-Note gets the bug fixed, he can use our test case to test it directly.'
$ws.Cells.Item(29, 7).Value = 'I always know I have a lot to learn. But through the course, I have some clues about what to learn currently. This course indeed needs team work in both thinking and doing, which made me understand how important it is. And I learned some skills about how to do better teamwork. All the things we learned in the course I will keep reviewing them once in a while. I believe I could get some new points of view when I gain more experience.'
$ws.Rows.Item(29).RowHeight = 247

# Update view state to match target selection/zoom
$excel.ActiveWindow.Zoom = 85
$ws.Range("G29").Select() | Out-Null
